$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (Förändrad / "Changed") holds a date serial value that was
# bumped by one day (46074 -> 46075) for every data row (rows 2-363).
$ws.Range("C2:C363").Value = 46075
